$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 61, shifting existing rows 61-164 down to 62-165
$ws.Rows.Item(61).Insert()

# Populate the new row 61 with the new weekly price-report entry
$ws.Cells.Item(61, 1).Value = 3
$ws.Cells.Item(61, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(61, 3).Value = "Coquimbo"
$ws.Cells.Item(61, 4).Value = 44645
$ws.Cells.Item(61, 5).Value = 5
$ws.Cells.Item(61, 6).Value = 100112052
$ws.Cells.Item(61, 7).Value = "Albahaca"
$ws.Cells.Item(61, 8).Value = "Sin especificar"
$ws.Cells.Item(61, 9).Value = "Primera"
$ws.Cells.Item(61, 10).Value = 60
$ws.Cells.Item(61, 11).Value = 4500
$ws.Cells.Item(61, 12).Value = 4500
$ws.Cells.Item(61, 13).Value = 4500
$ws.Cells.Item(61, 14).Value = "$/docena de matas"
$ws.Cells.Item(61, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(61, 16).Value = 750
$ws.Cells.Item(61, 17).Value = 6
$ws.Cells.Item(61, 18).Value = "Hortaliza"
